# Fixed capital letter error: add a "case" column (grammatical case used
# when rendering the certificate text) with the correct lower-case value
# "дательный" (dative case) for the first record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in N1 and corrected value in N2
$ws.Range("N1").Value = "case"
$ws.Range("N2").Value = "дательный"

# Reflect the view state captured after the edit: zoomed to 70% and the
# new cell N2 selected (instead of the old K1:L1 selection at 100%).
$excel.ActiveWindow.Zoom = 70
$ws.Range("N2").Select()
